$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-value updates derived from the authoritative diff, keyed by row -> column -> new text.
$rowUpdates = @(
    @{ Row = 2; D = '34.057.42'; E = '  -1.87%  ' },
    @{ Row = 3; D = '1.789.57'; E = '  -0.18%  ' },
    @{ Row = 4; D = '0.998'; E = '  -0.04%  ' },
    @{ Row = 5; D = '221.38'; E = '  -0.94%  ' },
    @{ Row = 6; D = '0.553'; E = '  +0.05%  ' },
    @{ Row = 7; E = '  -0.04%  ' },
    @{ Row = 8; D = '32.67'; E = '  +0.71%  ' },
    @{ Row = 9; E = '  +0.59%  ' },
    @{ Row = 10; E = '  +0.43%  ' },
    @{ Row = 11; E = '  -0.81%  ' },
    @{ Row = 12; D = '2.044.38'; E = '  -0.24%  ' },
    @{ Row = 13; D = '1.787.26'; E = '  -0.43%  ' },
    @{ Row = 14; D = '10.89'; E = '  -0.89%  ' },
    @{ Row = 15; D = '0.626'; E = '  -0.99%  ' },
    @{ Row = 16; D = '34.043.77'; E = '  -1.91%  ' },
    @{ Row = 17; E = '  -3.01%  ' },
    @{ Row = 18; D = '67.90'; E = '  -1.61%  ' },
    @{ Row = 19; D = '244.20'; E = '  -3.79%  ' },
    @{ Row = 20; E = '  -2.60%  ' },
    @{ Row = 21; E = '  +0.00%  ' },
    @{ Row = 22; D = '10.81'; E = '  +0.99%  ' },
    @{ Row = 23; E = '  -2.47%  ' },
    @{ Row = 24; E = '  -1.40%  ' },
    @{ Row = 25; D = '157.92'; E = '  -1.47%  ' },
    @{ Row = 26; D = '16.37'; E = '  +0.18%  ' },
    @{ Row = 27; D = '7.05'; E = '  -0.78%  ' },
    @{ Row = 28; E = '  -1.59%  ' },
    @{ Row = 29; E = '  -0.01%  ' },
    @{ Row = 30; D = '0.0520'; E = '  -1.44%  ' },
    @{ Row = 31; E = '  +0.83%  ' },
    @{ Row = 32; D = '3.67'; E = '  -3.07%  ' },
    @{ Row = 33; E = '  -2.97%  ' },
    @{ Row = 34; D = '1.81'; E = '  -2.81%  ' },
    @{ Row = 35; D = '1.395.41'; E = '  -2.77%  ' },
    @{ Row = 36; D = '0.640'; E = '  +0.86%  ' },
    @{ Row = 37; E = '  -0.42%  ' },
    @{ Row = 38; E = '  -3.35%  ' },
    @{ Row = 39; D = '79.65'; E = '  -5.90%  ' },
    @{ Row = 40; D = '0.922'; E = '  -0.30%  ' },
    @{ Row = 41; E = '  +0.94%  ' },
    @{ Row = 42; E = '  -2.64%  ' },
    @{ Row = 43; E = '  +1.16%  ' },
    @{ Row = 44; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '107.47'; E = '  +1.93%  ' },
    @{ Row = 45; B = 'WEMIXToken'; C = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; D = '1.05'; E = '  -1.01%  ' },
    @{ Row = 46; E = '  -1.25%  ' },
    @{ Row = 47; D = '0.0493'; E = '  +0.55%  ' },
    @{ Row = 48; D = '1.945.74'; E = '  +0.17%  ' },
    @{ Row = 49; B = 'InjectiveProtocol'; C = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D = '11.97'; E = '  +0.16%  ' },
    @{ Row = 50; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '0.998'; E = '  -0.09%  ' },
    @{ Row = 51; D = '0.0₆0128'; E = '  +2.02%  ' }
)

# A plain (unsigned, un-prefixed) decimal number written through COM gets
# auto-converted to a numeric cell by Excel, but the source data keeps these
# single-decimal price figures as literal text. Detect that shape and write
# it with a leading quote (forces text entry) then strip the resulting
# "quote prefix" cell style back to Normal so no stray formatting is left
# behind - only the cell's stored value/type should change, matching the diff.
function Test-PlainNumericText($s) {
    return $s -match '^[+-]?\d+(\.\d+)?$'
}

foreach ($update in $rowUpdates) {
    $rowNum = $update.Row
    foreach ($col in @("B", "C", "D", "E")) {
        if ($update.ContainsKey($col)) {
            $newValue = $update[$col]
            $cell = $ws.Range("$col$rowNum")
            if (Test-PlainNumericText $newValue) {
                $cell.Value = "'" + $newValue
                $cell.Style = "Normal"
            } else {
                $cell.Value = $newValue
            }
        }
    }
}
